$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "2021" year column header (R4) — same formatting as the existing
# P4/Q4 year headers.
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("R4").Value = 2021

# New data point (R5) — same formatting as the neighbouring Q5 cell, but
# with a dedicated "0.0" number format.
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("R5").Value = 102.20441221981518
$ws.Range("R5").NumberFormat = "0.0"

# Move the active selection, matching the saved view state.
$ws.Range("S9").Select() | Out-Null
